# Auto update stock data
# Updates the "Date_1" (column A) and "EBITDA" (column B) values for the
# most-recent-observation row of each company block, moving the reporting
# date from 2026/01/16 to 2026/01/17 and refreshing the EBITDA figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new EBITDA value (column B). Row 38 keeps its original EBITDA.
$updates = @{
    2  = "7.50"
    8  = "8.84"
    14 = "3.21"
    20 = "14.10"
    26 = "11.96"
    32 = "29.18"
    38 = $null
    44 = "17.13"
    50 = "12.75"
    56 = "32.82"
    62 = "12.41"
    68 = "13.46"
    74 = "19.76"
}

foreach ($row in $updates.Keys) {
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026/01/17"

    $ebitda = $updates[$row]
    if ($ebitda -ne $null) {
        $ebitdaCell = $ws.Cells.Item($row, 2)
        $ebitdaCell.NumberFormat = "@"
        $ebitdaCell.Value = $ebitda
    }
}
